$wb = $excel.ActiveWorkbook

# Update the contact name on the "Users" sheet
$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("A2").Value = "Ayati Arvind"

# Update the active selection on the "Users" sheet
$usersSheet.Activate()
$usersSheet.Range("D8").Select()
